# allocation.xlsx restructuring:
#  - "Dynamic" sheet renamed, relabeled bank ranges in hex, becomes active tab
#  - "Bank3" renamed to "Bank0x3", gains a "Code Purpose" column
#  - 12 new per-bank sheets (Bank0x1, Bank0x2, Bank0x4..Bank0xD) inserted
#  - BANK61 / BANK60 / Golden keep their data, just shift position to the end

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Dynamic" -> "Dynamic (Bank 0x10 - 0x27)" : relabel the bank-range
#    column from decimal ranges to hex ranges, make it the active tab.
# ---------------------------------------------------------------------
$dyn = $wb.Worksheets.Item("Dynamic")
$dyn.Name = "Dynamic (Bank 0x10 - 0x27)"
$dyn.Range("E2").Value = "0x10"
$dyn.Range("E3").Value = "0x11"
$dyn.Range("E4").Value = "0x12-0x18"
$dyn.Range("E5").Value = "0x19-0x21"
$dyn.Range("E6").Value = "0x22-0x27"

# ---------------------------------------------------------------------
# 2. "Bank3" -> "Bank0x3" : add a "Code Purpose" column, E2 becomes a
#    live formula instead of a stale literal.
# ---------------------------------------------------------------------
$bank3 = $wb.Worksheets.Item("Bank3")
$bank3.Name = "Bank0x3"
$bank3.Range("F1").Value = "Code Purpose"
$bank3.Range("E2").Formula = "=C2"
$bank3.Range("F2").Value = "commands.c commands"
$bank3.Columns.Item(6).AutoFit()

# ---------------------------------------------------------------------
# 3. Build the 12 new per-bank sheets from a template identical to the
#    "Bank0x3" layout (header row + one data row + totals row),
#    created in the order that matches their final sheetId allocation.
# ---------------------------------------------------------------------
function New-BankSheet {
    param([string]$Name, [double]$Size, [bool]$FormulaTotal, [double]$TotalValue, [string]$Purpose)

    $template = $wb.Worksheets.Item($wb.Worksheets.Item("Bank0x3").Index)
    $template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
    $new = $wb.Worksheets.Item($wb.Worksheets.Count)
    $new.Name = $Name

    # Clear the extra rows (3:5) that only Bank0x3 needs.
    $new.Range("A3:F5").ClearContents()

    $new.Range("C2").Value = $Size
    if ($FormulaTotal) {
        $new.Range("E2").Formula = "=C2"
    } else {
        $new.Range("E2").Value = $TotalValue
    }
    $new.Range("F2").Value = $Purpose
    $new.Columns.Item(6).AutoFit()
    $new.Range("A1").Select()

    return $new
}

New-BankSheet "Bank0x1" 3586 $true  3586 "View Code"                   | Out-Null
New-BankSheet "Bank0x2" 5123 $true  5123 "View Code"                   | Out-Null
New-BankSheet "Bank0x4" 3723 $true  3723 "View Code"                   | Out-Null
New-BankSheet "Bank0x5" 5123 $false 3093 "View Code"                   | Out-Null
New-BankSheet "Bank0x6" 6127 $true  6127 "View Code"                   | Out-Null
New-BankSheet "Bank0x7" 7554 $true  7554 "commands.c commands"         | Out-Null
New-BankSheet "Bank0x8" 7663 $true  7663 "commands.c commands"         | Out-Null
New-BankSheet "Bank0x9" 4234 $true  4234 "commands.c commands"         | Out-Null
New-BankSheet "Bank0xA" 6696 $true  6696 "commands.c commands"         | Out-Null
New-BankSheet "Bank0xB" 1153 $true  1153 "agiFiles.c Load Directories" | Out-Null
New-BankSheet "Bank0xC" 1645 $true  1645 "MEKA Main Logic"             | Out-Null
New-BankSheet "Bank0xD" 516  $true  516  "Logic Loading Logic.c"       | Out-Null

# ---------------------------------------------------------------------
# 4. Re-order every tab into its final position:
#    Dynamic, Bank0x1..Bank0xD (0x1,0x2,0x3,0x4..0xD), BANK61, BANK60, Golden
# ---------------------------------------------------------------------
$order = @("Bank0x1","Bank0x2","Bank0x3","Bank0x4","Bank0x5","Bank0x6","Bank0x7","Bank0x8","Bank0x9","Bank0xA","Bank0xB","Bank0xC","Bank0xD")
foreach ($nm in $order) {
    $wb.Worksheets.Item($nm).Move($wb.Worksheets.Item("BANK61"))
}

# ---------------------------------------------------------------------
# 5. Make the renamed "Dynamic" sheet the active tab with the selection
#    used after the edit.
# ---------------------------------------------------------------------
$dyn.Activate()
$dyn.Range("K22:K24").Select()
